# Reorders the team-name lists stored as text in column J ("new_top_teams")
# of Sheet1. Each distinct old list value is replaced by its corresponding
# new (reordered) list value, wherever it occurs in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$map = @{
    "['France', 'Belgium']" = "['Belgium', 'France']"
    "['Spain', 'Portugal']" = "['Portugal', 'Spain']"
    "['Ireland', 'Soviet Union']" = "['Soviet Union', 'Ireland']"
    "['Sweden', 'France']" = "['France', 'Sweden']"
    "['Sweden', 'England']" = "['England', 'Sweden']"
    "['Scotland', 'England']" = "['England', 'Scotland']"
    "['Spain', 'France']" = "['France', 'Spain']"
    "['Romania', 'Portugal']" = "['Portugal', 'Romania']"
    "['Italy', 'Belgium']" = "['Belgium', 'Italy']"
    "['Spain', 'Yugoslavia']" = "['Yugoslavia', 'Spain']"
    "['Spain', 'Greece']" = "['Greece', 'Spain']"
    "['France', 'England']" = "['England', 'France']"
    "['Czech Republic', 'Netherlands']" = "['Netherlands', 'Czech Republic']"
    "['Turkey', 'Portugal']" = "['Portugal', 'Turkey']"
    "['Russia', 'Czech Republic']" = "['Czech Republic', 'Russia']"
    "['Russia', 'Greece']" = "['Greece', 'Russia']"
    "['Germany', 'Portugal']" = "['Portugal', 'Germany']"
    "['Spain', 'Croatia']" = "['Croatia', 'Spain']"
    "['Spain', 'Italy']" = "['Italy', 'Spain']"
    "['Romania', 'France', 'Switzerland']" = "['Switzerland', 'Romania', 'France']"
    "['Albania', 'France', 'Switzerland']" = "['Switzerland', 'France', 'Albania']"
    "['Wales', 'Slovakia', 'England']" = "['Wales', 'England', 'Slovakia']"
    "['Germany', 'Poland', 'Northern Ireland']" = "['Germany', 'Northern Ireland', 'Poland']"
    "['Spain', 'Czech Republic', 'Croatia']" = "['Croatia', 'Spain', 'Czech Republic']"
    "['Spain', 'Turkey', 'Croatia']" = "['Turkey', 'Croatia', 'Spain']"
    "['Sweden', 'Italy', 'Belgium']" = "['Belgium', 'Italy', 'Sweden']"
    "['Italy', 'Belgium', 'Ireland']" = "['Belgium', 'Ireland', 'Italy']"
    "['Portugal', 'Hungary', 'Iceland']" = "['Portugal', 'Iceland', 'Hungary']"
    "['Italy', 'Wales', 'Switzerland']" = "['Wales', 'Italy', 'Switzerland']"
    "['Russia', 'Finland', 'Belgium']" = "['Belgium', 'Russia', 'Finland']"
    "['Finland', 'Belgium', 'Denmark']" = "['Belgium', 'Denmark', 'Finland']"
    "['Russia', 'Belgium', 'Denmark']" = "['Belgium', 'Denmark', 'Russia']"
    "['Ukraine', 'Austria', 'Netherlands']" = "['Ukraine', 'Netherlands', 'Austria']"
    "['Czech Republic', 'Croatia', 'England']" = "['Croatia', 'England', 'Czech Republic']"
    "['Sweden', 'Spain', 'Slovakia']" = "['Slovakia', 'Sweden', 'Spain']"
    "['Spain', 'Slovakia', 'Sweden']" = "['Slovakia', 'Spain', 'Sweden']"
    "['Germany', 'Portugal', 'France']" = "['Portugal', 'Germany', 'France']"
    "['Portugal', 'Hungary', 'France']" = "['Portugal', 'France', 'Hungary']"
    "['Germany', 'Scotland', 'Switzerland']" = "['Switzerland', 'Germany', 'Scotland']"
    "['Germany', 'Hungary', 'Switzerland']" = "['Switzerland', 'Germany', 'Hungary']"
    "['Spain', 'Italy', 'Albania']" = "['Italy', 'Spain', 'Albania']"
    "['Spain', 'Italy', 'Croatia']" = "['Croatia', 'Italy', 'Spain']"
    "['Denmark', 'Slovenia', 'England']" = "['Slovenia', 'England', 'Denmark']"
    "['Austria', 'France', 'Netherlands']" = "['France', 'Netherlands', 'Austria']"
    "['Netherlands', 'France', 'Austria']" = "['France', 'Netherlands', 'Austria']"
    "['Romania', 'Belgium', 'Slovakia']" = "['Belgium', 'Romania', 'Slovakia']"
    "['Ukraine', 'Slovakia', 'Belgium']" = "['Ukraine', 'Belgium', 'Slovakia']"
    "['Turkey', 'Portugal', 'Czech Republic']" = "['Portugal', 'Turkey', 'Czech Republic']"
    "['Turkey', 'Portugal', 'Georgia']" = "['Portugal', 'Turkey', 'Georgia']"
    "['Georgia', 'Portugal', 'Turkey']" = "['Portugal', 'Turkey', 'Georgia']"
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 10  # column J = "new_top_teams"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
